# Edit script: add age_group_order_* columns (E,F,G) to the age_groups sheet,
# remove tabSelected from lsh_sheet_names sheet, make age_groups the active/selected
# sheet with a new selection, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("age_groups")

# Header row (row 1): new column headers in E1:G1
$ws.Range("E1").Value = "age_group_order_official"
$ws.Range("F1").Value = "age_group_order_three"
$ws.Range("G1").Value = "age_group_order_simple"

# Data rows 2..122 correspond to age values 0..120 in column A.
# E = ordinal of the official decade-based age group (column B)
# F = ordinal of the three-way age group (column C): 0-50 / 51-74 / 75+
# G = ordinal of the simple two-way age group (column D): 0-50 / 51+
$data = @(
    @(1,1,1),
    @(1,1,1),
    @(1,1,1),
    @(1,1,1),
    @(1,1,1),
    @(1,1,1),
    @(1,1,1),
    @(1,1,1),
    @(1,1,1),
    @(1,1,1),
    @(2,1,1),
    @(2,1,1),
    @(2,1,1),
    @(2,1,1),
    @(2,1,1),
    @(2,1,1),
    @(2,1,1),
    @(2,1,1),
    @(2,1,1),
    @(2,1,1),
    @(3,1,1),
    @(3,1,1),
    @(3,1,1),
    @(3,1,1),
    @(3,1,1),
    @(3,1,1),
    @(3,1,1),
    @(3,1,1),
    @(3,1,1),
    @(3,1,1),
    @(4,1,1),
    @(4,1,1),
    @(4,1,1),
    @(4,1,1),
    @(4,1,1),
    @(4,1,1),
    @(4,1,1),
    @(4,1,1),
    @(4,1,1),
    @(4,1,1),
    @(5,1,1),
    @(5,1,1),
    @(5,1,1),
    @(5,1,1),
    @(5,1,1),
    @(5,1,1),
    @(5,1,1),
    @(5,1,1),
    @(5,1,1),
    @(5,1,1),
    @(6,1,1),
    @(6,2,2),
    @(6,2,2),
    @(6,2,2),
    @(6,2,2),
    @(6,2,2),
    @(6,2,2),
    @(6,2,2),
    @(6,2,2),
    @(6,2,2),
    @(7,2,2),
    @(7,2,2),
    @(7,2,2),
    @(7,2,2),
    @(7,2,2),
    @(7,2,2),
    @(7,2,2),
    @(7,2,2),
    @(7,2,2),
    @(7,2,2),
    @(8,2,2),
    @(8,2,2),
    @(8,2,2),
    @(8,2,2),
    @(8,2,2),
    @(8,3,2),
    @(8,3,2),
    @(8,3,2),
    @(8,3,2),
    @(8,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2),
    @(9,3,2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Range("E$row").Value = $vals[0]
    $ws.Range("F$row").Value = $vals[1]
    $ws.Range("G$row").Value = $vals[2]
}

# Remove the active-tab marker from lsh_sheet_names and move it to age_groups,
# and update the selected cell there to I8.
$namesWs = $wb.Worksheets.Item("lsh_sheet_names")
$namesWs.Activate()

$ws.Activate()
$ws.Range("I8").Select()

Write-Host "done"
